$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J7").Value = "slightly lower"
$ws.Range("J8").Value = "lower"
$ws.Range("J9").Value = "much lower"
$ws.Range("J10").Value = "[Manually written]"

$ws.Range("M10").Select()
